$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 61-77 (A:E) of the DSTU2/R4 table are being restructured:
#  - two new "US Core ADI" related rows are inserted right after the
#    "Treatment Intervention Preference Profile" row (old row 60),
#    pushing the "**Extensions**" section heading and everything below
#    it down by two rows;
#  - two new extension rows ("US Core Authentication Time Extension" and
#    "US Core Interpreter Required Extension") are appended near the end
#    of the Extensions section, just before the closing footnote row.
# The net result is that what used to be rows 61-73 now occupies rows
# 63-77 (shifted by 2), and a new block occupies rows 61-62, with two
# additional rows (75 and 76) inserted before the final footnote row.

$rows = @(
    ,@(61, "     -", "US Core Observation ADI Documentation Profile", "8.0.0", "", "")
    ,@(62, "     -", "US Core ADI DocumentReference Profile", "8.0.0", "", "")
    ,@(63, "**Extensions**", "", "", "", "")
    ,@(64, "Sex of patient assigned at birth", "US Core Birth Sex Extension", "3.1.1*", "", "After version 6.0.0 this extension is no longer a USCDI requirement.")
    ,@(65, "Argonaut ethnicity Extension", "US Core Ethnicity Extension", "3.1.1*", "", "")
    ,@(66, "Argonaut Race Extension", "US Core Race Extension", "3.1.1*", "", "")
    ,@(67, "     -", "US Core Direct email Extension", "3.1.1*", "", "")
    ,@(68, "     -", "US Core Extension Questionnaire URI", "5.0.0", "", "")
    ,@(69, "     -", "US Core Gender Identity Extension", "5.0.0", "", "")
    ,@(70, "     -", "US Core Tribal Affiliation Extension", "6.0.0", "", "")
    ,@(71, "     -", "US Core Jurisdiction Extension", "6.0.0", "", "")
    ,@(72, "     -", "US Core USCDI Requirements Extension", "6.0.0", "", "This extension is only used on US Core Profile StructureDefinition elements")
    ,@(73, "     -", "US Core Sex Extension", "6.1.0", "", "")
    ,@(74, "     -", "US Core Medication Adherence Extension", "7.0.0", "", "")
    ,@(75, "     ", "US Core Authentication Time Extension", "8.0.0", "", "")
    ,@(76, "     -", "US Core Interpreter Required Extension", "8.0.0", "", "")
    ,@(77, "", "", "", "", "\* 3.1.1  *or prior* version of US Core")
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($i = 1; $i -lt $row.Length; $i++) {
        $col = $i
        $ws.Cells.Item($r, $col).Value = $row[$i]
    }
}
